$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 2 (this shifts all existing data rows
# 2..10 down to 3..11, which reproduces the rest of the diff automatically).
$ws.Rows.Item(2).Insert()

# The inserted row inherited the bold header formatting from row 1 (the
# interop layer always copies formats from the row above on Insert).
# Strip that back to "no style" to match the rest of the data rows...
$ws.Range("A2:R2").ClearFormats()

# ...then restore just the date-column number format on D2 by copying it
# from the date cell directly below (which already carries style index 2).
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's contents.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44473
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 9500
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9750
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 390
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
